$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.859.04"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "2.624.70"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.550"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("D9").Value = "2.624.19"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  +9.15%  "
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.347"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("E15").Value = "  +3.66%  "
$ws.Range("D16").Value = "3.105.24"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "67.820.47"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "2.642.93"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "370.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.11%  "
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("E29").Value = "  +2.31%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "576.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.94%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.72%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").Value = "0.0₆0333"
$ws.Range("E43").Value = "  +15.04%  "
$ws.Range("E44").Value = "  +2.93%  "
$ws.Range("E45").Value = "  +6.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "155.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("E51").Value = "  -1.26%  "
